$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bulk-updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (2-171).
for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
